$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 48
$ws.Range("E10").Value = 27
$ws.Range("D11").Value = 55
$ws.Range("E11").Value = 34
$ws.Range("D12").Value = 59
$ws.Range("E12").Value = 38
$ws.Range("D13").Value = 66
$ws.Range("E13").Value = 45
$ws.Range("D14").Value = 68
$ws.Range("E14").Value = 33
$ws.Range("D15").Value = 82
$ws.Range("E15").Value = 36
$ws.Range("F15").Value = 25
$ws.Range("D16").Value = 105
$ws.Range("E16").Value = 51
$ws.Range("F16").Value = 33
$ws.Range("D17").Value = 118
$ws.Range("E17").Value = 59
$ws.Range("D18").Value = 140
$ws.Range("E18").Value = 74
$ws.Range("F18").Value = 45
$ws.Range("D19").Value = 168
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 33
$ws.Range("C20").Value = 46
$ws.Range("D20").Value = 195
$ws.Range("E20").Value = 118
$ws.Range("F20").Value = 31
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 231
$ws.Range("E21").Value = 134
$ws.Range("F21").Value = 43
$ws.Range("D22").Value = 291
$ws.Range("E22").Value = 178
$ws.Range("F22").Value = 54
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 66
$ws.Range("D23").Value = 351
$ws.Range("E23").Value = 222
$ws.Range("F23").Value = 63
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = 68
$ws.Range("D24").Value = 425
$ws.Range("E24").Value = 267
$ws.Range("F24").Value = 90
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 77
$ws.Range("D25").Value = 538
$ws.Range("E25").Value = 356
$ws.Range("F25").Value = 105
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = 97
$ws.Range("D26").Value = 671
$ws.Range("E26").Value = 454
$ws.Range("F26").Value = 120
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 113
$ws.Range("D27").Value = 799
$ws.Range("E27").Value = 534
$ws.Range("F27").Value = 152
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = 129
$ws.Range("D28").Value = 949
$ws.Range("E28").Value = 630
$ws.Range("F28").Value = 190
$ws.Range("B29").Value = 4
$ws.Range("C29").Value = 158
$ws.Range("D29").Value = 1142
$ws.Range("E29").Value = 751
$ws.Range("F29").Value = 233
$ws.Range("B30").Value = 5
$ws.Range("C30").Value = 182
$ws.Range("D30").Value = 1352
$ws.Range("E30").Value = 863
$ws.Range("F30").Value = 307
$ws.Range("B31").Value = 5
$ws.Range("C31").Value = 217
$ws.Range("D31").Value = 1610
$ws.Range("E31").Value = 999
$ws.Range("F31").Value = 394
$ws.Range("B32").Value = 6
$ws.Range("C32").Value = 265
$ws.Range("D32").Value = 1901
$ws.Range("E32").Value = 1168
$ws.Range("F32").Value = 468
$ws.Range("B33").Value = 7
$ws.Range("C33").Value = 319
$ws.Range("D33").Value = 2271
$ws.Range("E33").Value = 1392
$ws.Range("F33").Value = 560
$ws.Range("B34").Value = 7
$ws.Range("C34").Value = 391
$ws.Range("D34").Value = 2707
$ws.Range("E34").Value = 1644
$ws.Range("F34").Value = 672
$ws.Range("B35").Value = 9
$ws.Range("C35").Value = 489
$ws.Range("D35").Value = 3231
$ws.Range("E35").Value = 1976
$ws.Range("F35").Value = 766
$ws.Range("B36").Value = 11
$ws.Range("C36").Value = 611
$ws.Range("D36").Value = 3769
$ws.Range("E36").Value = 2277
$ws.Range("F36").Value = 881
$ws.Range("B37").Value = 13
$ws.Range("C37").Value = 733
$ws.Range("D37").Value = 4422
$ws.Range("E37").Value = 2653
$ws.Range("F37").Value = 1036
$ws.Range("B38").Value = 15
$ws.Range("C38").Value = 879
$ws.Range("D38").Value = 5203
$ws.Range("E38").Value = 3111
$ws.Range("F38").Value = 1213
$ws.Range("B39").Value = 19
$ws.Range("C39").Value = 1063
$ws.Range("D39").Value = 5968
$ws.Range("E39").Value = 3474
$ws.Range("F39").Value = 1431
$ws.Range("B40").Value = 22
$ws.Range("C40").Value = 1255
$ws.Range("D40").Value = 6954
$ws.Range("E40").Value = 3979
$ws.Range("F40").Value = 1720
$ws.Range("B41").Value = 25
$ws.Range("C41").Value = 1492
$ws.Range("D41").Value = 8196
$ws.Range("E41").Value = 4687
$ws.Range("F41").Value = 2017
$ws.Range("B42").Value = 34
$ws.Range("C42").Value = 1769
$ws.Range("D42").Value = 9606
$ws.Range("E42").Value = 5500
$ws.Range("F42").Value = 2337
$ws.Range("B43").Value = 42
$ws.Range("C43").Value = 2092
$ws.Range("D43").Value = 11165
$ws.Range("E43").Value = 6365
$ws.Range("F43").Value = 2708
$ws.Range("B44").Value = 52
$ws.Range("C44").Value = 2494
$ws.Range("D44").Value = 12887
$ws.Range("E44").Value = 7269
$ws.Range("F44").Value = 3124
$ws.Range("B45").Value = 72
$ws.Range("C45").Value = 2975
$ws.Range("D45").Value = 14876
$ws.Range("E45").Value = 8376
$ws.Range("F45").Value = 3525
